$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume values per diff
$ws.Range("D2").Value = '37.928.70'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '2.038.74'
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''228.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("E6").Value = '  -0.65%  '
$ws.Range("D7").Value = '''60.88'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.48%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '''0.379'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.67%  '
$ws.Range("E10").Value = '  +0.79%  '
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("D12").Value = '2.337.95'
$ws.Range("E12").Value = '  -0.75%  '
$ws.Range("D13").Value = '''14.54'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").Value = '''21.49'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.61%  '
$ws.Range("E15").Value = '  +1.52%  '
$ws.Range("D16").Value = '''5.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.82%  '
$ws.Range("D17").Value = '2.028.66'
$ws.Range("E17").Value = '  -1.17%  '
$ws.Range("D18").Value = '37.866.45'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").Value = '''69.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").Value = '''5.91'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.42%  '
$ws.Range("E21").Value = '  -1.38%  '
$ws.Range("D22").Value = '''224.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").Value = '''2.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("D26").Value = '''9.35'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("D27").Value = '''167.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.79%  '
$ws.Range("E28").Value = '  -1.95%  '
$ws.Range("D29").Value = '''18.91'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.54%  '
$ws.Range("E30").Value = '  -3.16%  '
$ws.Range("E31").Value = '  +0.80%  '
$ws.Range("D32").Value = '''2.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.73%  '
$ws.Range("E33").Value = '  -2.62%  '
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("E35").Value = '  -1.53%  '
$ws.Range("D36").Value = '''6.39'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.08%  '
$ws.Range("D37").Value = '''2.30'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.63%  '
$ws.Range("D38").Value = '''3.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.88%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").Value = '1.540.89'
$ws.Range("E40").Value = '  +0.32%  '
$ws.Range("D41").Value = '''17.77'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.96%  '
$ws.Range("D42").Value = '''0.0219'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.94%  '
$ws.Range("D43").Value = '''96.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.56%  '
$ws.Range("D44").Value = '''2.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.82%  '
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("E46").Value = '  -2.21%  '
$ws.Range("D47").Value = '''4.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.17%  '
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("D50").Value = '''7.10'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").Value = '2.227.41'
$ws.Range("E51").Value = '  -0.65%  '
